$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format before writing numeric-looking strings so Excel
# does not silently coerce them into floating point cells, then restore the
# default "Normal" style so no stray style index is left on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '58.724.35'
Set-TextValue $ws.Range('E2') '  -2.05%  '
Set-TextValue $ws.Range('D3') '2.304.45'
Set-TextValue $ws.Range('E3') '  -4.13%  '
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '547.16'
Set-TextValue $ws.Range('E5') '  -1.15%  '
Set-TextValue $ws.Range('E6') '  -2.45%  '
Set-TextValue $ws.Range('E7') '  -0.01%  '
Set-TextValue $ws.Range('D8') '0.572'
Set-TextValue $ws.Range('E8') '  -1.79%  '
Set-TextValue $ws.Range('D9') '2.302.66'
Set-TextValue $ws.Range('E9') '  -4.16%  '
Set-TextValue $ws.Range('E10') '  -2.79%  '
Set-TextValue $ws.Range('E11') '  -1.83%  '
Set-TextValue $ws.Range('E12') '  +1.41%  '
Set-TextValue $ws.Range('E13') '  -4.46%  '
Set-TextValue $ws.Range('D14') '23.91'
Set-TextValue $ws.Range('E14') '  -2.66%  '
Set-TextValue $ws.Range('D15') '2.716.17'
Set-TextValue $ws.Range('E15') '  -4.23%  '
Set-TextValue $ws.Range('D16') '58.711.69'
Set-TextValue $ws.Range('E16') '  -1.83%  '
Set-TextValue $ws.Range('E17') '  -2.89%  '
Set-TextValue $ws.Range('D18') '2.312.22'
Set-TextValue $ws.Range('E18') '  -3.73%  '
Set-TextValue $ws.Range('D19') '10.69'
Set-TextValue $ws.Range('E19') '  -4.05%  '
Set-TextValue $ws.Range('E20') '  -3.84%  '
Set-TextValue $ws.Range('D21') '314.73'
Set-TextValue $ws.Range('E21') '  -3.24%  '
Set-TextValue $ws.Range('D22') '6.46'
Set-TextValue $ws.Range('E22') '  -4.15%  '
Set-TextValue $ws.Range('E23') '  -0.01%  '
Set-TextValue $ws.Range('E24') '  -1.83%  '
Set-TextValue $ws.Range('E25') '  -6.32%  '
Set-TextValue $ws.Range('E26') '  +0.07%  '
Set-TextValue $ws.Range('D27') '8.11'
Set-TextValue $ws.Range('E27') '  -5.40%  '
Set-TextValue $ws.Range('E28') '  -5.14%  '
Set-TextValue $ws.Range('E29') '  -1.62%  '
Set-TextValue $ws.Range('D30') '168.14'
Set-TextValue $ws.Range('E30') '  -1.01%  '
Set-TextValue $ws.Range('E31') '  -4.78%  '
Set-TextValue $ws.Range('E32') '  +1.00%  '
Set-TextValue $ws.Range('D33') '5.79'
Set-TextValue $ws.Range('E33') '  -5.21%  '
Set-TextValue $ws.Range('E34') '  -4.55%  '
Set-TextValue $ws.Range('E36') '  -3.11%  '
Set-TextValue $ws.Range('E37') '  -0.03%  '
Set-TextValue $ws.Range('E38') '  -4.12%  '
Set-TextValue $ws.Range('E39') '  -4.85%  '
Set-TextValue $ws.Range('D40') '38.10'
Set-TextValue $ws.Range('E40') '  -1.27%  '
Set-TextValue $ws.Range('E41') '  -4.68%  '
Set-TextValue $ws.Range('D42') '298.10'
Set-TextValue $ws.Range('E42') '  -7.39%  '
Set-TextValue $ws.Range('D43') '141.57'
Set-TextValue $ws.Range('E43') '  -3.40%  '
Set-TextValue $ws.Range('E44') '  -3.66%  '
Set-TextValue $ws.Range('E45') '  -1.14%  '
Set-TextValue $ws.Range('D46') '0.0501'
Set-TextValue $ws.Range('E46') '  -2.25%  '
Set-TextValue $ws.Range('D47') '0.556'
Set-TextValue $ws.Range('E47') '  -3.06%  '
Set-TextValue $ws.Range('D48') '18.52'
Set-TextValue $ws.Range('E48') '  -6.27%  '
Set-TextValue $ws.Range('E49') '  -2.42%  '
Set-TextValue $ws.Range('D50') '16.64'
Set-TextValue $ws.Range('E50') '  -3.39%  '
Set-TextValue $ws.Range('D51') '11.01'
Set-TextValue $ws.Range('E51') '  -0.27%  '
